# DNS master/slave servicos publicos - renumber the LAN host addresses
# 10.0.0.4 -> 10.0.0.8, 10.0.0.5 -> 10.0.0.9, 10.0.0.6 -> 10.0.0.10, 10.0.0.7 -> 10.0.0.11
$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "4"; New = "8" },
    @{ Old = "5"; New = "9" },
    @{ Old = "6"; New = "10" },
    @{ Old = "7"; New = "11" }
)

foreach ($pair in $replacements) {
    $found = $d.Content.Find.Execute($pair.Old, $true, $true, $false, $false, $false,
                                      $true, 1, $false, $pair.New, 2)
    Write-Output ("Replaced '" + $pair.Old + "' -> '" + $pair.New + "': " + $found)
}
